# Append the "solved" continuation after the paragraph that ends with
# "...and get the answer as it correlates to the number." The new
# material adds the formula (8x + 1), works the example for counting to
# ten, and states the final answer (C).

$d = $word.ActiveDocument

# Locate the paragraph whose text ends with the "...correlates to the
# number." sentence -- this is where the bookmarkStart/bookmarkEnd
# ("_GoBack") live, and everything new gets appended right after it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "correlates to the number") {
        $anchorIndex = $i
    }
}

$cur = $d.Paragraphs.Item($anchorIndex).Range

# Blank paragraph
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 1).Range

# "The formula for this would be 8x +1 "
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 2).Range
$cur.InsertAfter("The formula for this would be 8x +1 ")

# Blank paragraph
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 3).Range

# "Based on this if the girl counted to ten the answer would be. "
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 4).Range
$cur.InsertAfter("Based on this if the girl counted to ten the answer would be. ")

# "10/8=1.25 1.25-.25=1 1*8=8 8+2=10 "
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 5).Range
$cur.InsertAfter("10/8=1.25 1.25-.25=1 1*8=8 8+2=10 ")

# "so since A = 8x +1 the answer would be C"
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 6).Range
$cur.InsertAfter("so since A = 8x +1 the answer would be C")

# Blank paragraph
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 7).Range

# "If the girl counts to ten she will end on C."
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($anchorIndex + 8).Range
$cur.InsertAfter("If the girl counts to ten she will end on C.")
